# Generate Report for Handback
# Adds a new "in sync with en-US" row for 82752071-1c07-4ffe-91f3-7a0b4128ebe5.md
# to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$guid = "82752071-1c07-4ffe-91f3-7a0b4128ebe5"
$mdName = "$guid.md"
$pathAndName = "e2e\$guid.md"

$hyperFontColor = 15570276   # RGB(0x64,0x95,0xED) == existing "HyperLink" font colour FF6495ED
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$rOverview = $loOverview.Range.Rows.Count + $loOverview.Range.Row - 1

$wsOverview.Range("A$rOverview").Value = $mdName
$wsOverview.Range("B$rOverview").Value = $pathAndName
$wsOverview.Range("B$rOverview").Font.Underline = 2
$wsOverview.Range("B$rOverview").Font.Color = $hyperFontColor
$wsOverview.Range("C$rOverview").Value = ".md"
$wsOverview.Range("E$rOverview").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F$rOverview").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G$rOverview").Value = "2016-08-13 10:55:21"
$wsOverview.Range("G$rOverview").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B$rOverview"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/master/e2e/$guid.md", "", "", $pathAndName) | Out-Null

# ---------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$rZh = $loZh.Range.Rows.Count + $loZh.Range.Row - 1

$wsZh.Range("A$rZh").Value = $mdName
$wsZh.Range("A$rZh").Font.Underline = 2
$wsZh.Range("A$rZh").Font.Color = $hyperFontColor
$wsZh.Range("B$rZh").Value = ".md"
$wsZh.Range("C$rZh").Value = "Handed back: in sync with en-US"
$wsZh.Range("D$rZh").Value = "e2e"
$wsZh.Range("E$rZh").Value = "ht"
$wsZh.Range("F$rZh").Value = "True"
$wsZh.Range("G$rZh").Value = "$guid.6e7fac94290f194893159599289946e9174f74b9.zh-cn.xlf"
$wsZh.Range("H$rZh").Value = "2016-08-13 10:55:13"
$wsZh.Range("H$rZh").NumberFormat = $dateFmt
$wsZh.Range("I$rZh").Value = $mdName
$wsZh.Range("I$rZh").Font.Underline = 2
$wsZh.Range("I$rZh").Font.Color = $hyperFontColor
$wsZh.Range("J$rZh").Value = "$guid.6e7fac94290f194893159599289946e9174f74b9.zh-cn.xlf"
$wsZh.Range("K$rZh").Value = "2016-08-13 10:55:42"
$wsZh.Range("K$rZh").NumberFormat = $dateFmt
$wsZh.Range("L$rZh").Value = ""
$wsZh.Range("M$rZh").Value = "True"
$wsZh.Range("N$rZh").Value = ""
$wsZh.Range("O$rZh").Value = "False"
$wsZh.Range("P$rZh").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A$rZh"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/master/e2e/$guid.md", "", "", $mdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I$rZh"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/$guid.md", "", "", $mdName) | Out-Null

# ---------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$rDe = $loDe.Range.Rows.Count + $loDe.Range.Row - 1

$wsDe.Range("A$rDe").Value = $mdName
$wsDe.Range("A$rDe").Font.Underline = 2
$wsDe.Range("A$rDe").Font.Color = $hyperFontColor
$wsDe.Range("B$rDe").Value = ".md"
$wsDe.Range("C$rDe").Value = "Handed back: in sync with en-US"
$wsDe.Range("D$rDe").Value = "e2e"
$wsDe.Range("E$rDe").Value = "ht"
$wsDe.Range("F$rDe").Value = "True"
$wsDe.Range("G$rDe").Value = "$guid.6e7fac94290f194893159599289946e9174f74b9.de-de.xlf"
$wsDe.Range("H$rDe").Value = "2016-08-13 10:55:21"
$wsDe.Range("H$rDe").NumberFormat = $dateFmt
$wsDe.Range("I$rDe").Value = $mdName
$wsDe.Range("I$rDe").Font.Underline = 2
$wsDe.Range("I$rDe").Font.Color = $hyperFontColor
$wsDe.Range("J$rDe").Value = "$guid.6e7fac94290f194893159599289946e9174f74b9.de-de.xlf"
$wsDe.Range("K$rDe").Value = "2016-08-13 10:55:52"
$wsDe.Range("K$rDe").NumberFormat = $dateFmt
$wsDe.Range("L$rDe").Value = ""
$wsDe.Range("M$rDe").Value = "True"
$wsDe.Range("N$rDe").Value = ""
$wsDe.Range("O$rDe").Value = "False"
$wsDe.Range("P$rDe").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A$rDe"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/master/e2e/$guid.md", "", "", $mdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I$rDe"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/$guid.md", "", "", $mdName) | Out-Null
